$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '49.403.33'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '2.625.66'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue 'D5' '111.23'
$ws.Range('E5').Value = '  -2.18%  '
Set-TextValue 'D6' '324.27'
$ws.Range('E6').Value = '  -1.01%  '
Set-TextValue 'D7' '0.521'
$ws.Range('E7').Value = '  -1.87%  '
Set-TextValue 'D8' '0.999'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -3.60%  '
Set-TextValue 'D10' '39.31'
$ws.Range('E10').Value = '  -5.11%  '
Set-TextValue 'D11' '20.03'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('E13').Value = '  +1.32%  '
Set-TextValue 'D14' '7.29'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '3.039.63'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '2.625.07'
$ws.Range('E16').Value = '  -2.24%  '
Set-TextValue 'D17' '0.848'
$ws.Range('E17').Value = '  -3.37%  '
$ws.Range('D18').Value = '49.346.15'
$ws.Range('E18').Value = '  -1.30%  '
Set-TextValue 'D19' '12.90'
$ws.Range('E19').Value = '  -3.30%  '
Set-TextValue 'D20' '2.90'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('E21').Value = '  -2.26%  '
$ws.Range('E22').Value = '  -2.85%  '
Set-TextValue 'D23' '267.89'
$ws.Range('E23').Value = '  -4.16%  '
$ws.Range('E24').Value = '  -5.77%  '
$ws.Range('E25').Value = '  -2.96%  '
$ws.Range('E26').Value = '  +0.05%  '
Set-TextValue 'D27' '25.91'
$ws.Range('E27').Value = '  -4.12%  '
Set-TextValue 'D28' '10.11'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('E30').Value = '  -4.56%  '
Set-TextValue 'D31' '34.39'
$ws.Range('E31').Value = '  -6.67%  '
Set-TextValue 'D32' '49.44'
$ws.Range('E32').Value = '  -2.06%  '
Set-TextValue 'D33' '5.48'
$ws.Range('E33').Value = '  +0.37%  '
Set-TextValue 'D34' '0.0807'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  -0.12%  '
Set-TextValue 'D36' '18.73'
$ws.Range('E36').Value = '  -5.18%  '
Set-TextValue 'D37' '4.89'
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('E39').Value = '  -0.31%  '
Set-TextValue 'D40' '128.10'
$ws.Range('E40').Value = '  +4.02%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D41' '22.47'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D42' '0.110'
$ws.Range('E42').Value = '  -2.15%  '
Set-TextValue 'D43' '0.0321'
$ws.Range('E43').Value = '  +1.67%  '
Set-TextValue 'D44' '2.13'
$ws.Range('E44').Value = '  -5.42%  '
$ws.Range('D45').Value = '2.040.26'
$ws.Range('E45').Value = '  -1.99%  '
Set-TextValue 'D46' '2.15'
$ws.Range('E46').Value = '  +7.50%  '
Set-TextValue 'D47' '3.19'
$ws.Range('E47').Value = '  -5.53%  '
$ws.Range('E48').Value = '  -4.11%  '
Set-TextValue 'D49' '8.83'
$ws.Range('E49').Value = '  -3.94%  '
$ws.Range('E51').Value = '  +0.38%  '
